# "add wiz note loader"
#
# Appends more sample rows to Sheet1 (A3:B21) and Sheet2 (A3:B10), each
# driven by a NOW() watermark formula in column B (mirroring the existing
# A2/B2 seed row on both sheets), then leaves the UI focused on Sheet1.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet1: rows 3..21 -> A = 124..142, B = "=NOW()" (fills as one shared
# formula group, same shape Excel itself would produce for a fill-down).
# ---------------------------------------------------------------------
$r = 3
for ($val = 124; $val -le 142; $val++) {
    $ws1.Cells.Item($r, 1).Value = $val
    $r++
}

# Pull B3:B21's number format from the existing B2 watermark cell first
# (so no new style gets minted), THEN write the formulas.
$ws1.Range("B2").Copy() | Out-Null
$ws1.Range("B3:B21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws1.Range("B3:B21").Formula = "=NOW()"

# ---------------------------------------------------------------------
# Sheet2: rows 3..10 -> A = 4790..4797; only rows 3 & 4 also get a B
# watermark formula, rows 5..10 are A-only.
# ---------------------------------------------------------------------
$ws2.Range("A3").Value = 4790
$ws2.Range("A4").Value = 4791
$ws2.Range("A5").Value = 4792
$ws2.Range("A6").Value = 4793
$ws2.Range("A7").Value = 4794
$ws2.Range("A8").Value = 4795
$ws2.Range("A9").Value = 4796
$ws2.Range("A10").Value = 4797

$ws2.Range("B2").Copy() | Out-Null
$ws2.Range("B3:B4").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$ws2.Range("B3").Formula = "=NOW()"
$ws2.Range("B4").Formula = "=NOW()"

# ---------------------------------------------------------------------
# Selection / active-sheet state: Sheet2 ends up parked on B4, Sheet1
# becomes the active tab with focus on B7 (last Select() wins the tab).
# ---------------------------------------------------------------------
$ws2.Range("B4").Select() | Out-Null
$ws1.Range("A2:B21").Select() | Out-Null
$ws1.Range("B7").Select() | Out-Null
